# Config cleanup based on the Branch Type finished.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Info")

# Row 7 (Design): switch the branch type from BASE to FLOW
$ws.Range("B7").Value() = "FLOW"

# Row 8 (Converged router): no longer a production router -> FALSE,
# and add an explanatory note that it's only relevant for NAM
$ws.Range("B8").Value() = $false
$ws.Range("D8").Value() = "Only for NAM - has no effect for the other regions"

# Row 9 (Migration from MPLS): switch from the free-text option to a boolean FALSE,
# and add the same explanatory note
$ws.Range("B9").Value() = $false
$ws.Range("D9").Value() = "Only for NAM - has no effect for the other regions"

# Row 10 (ZBFW): no longer enabled
$ws.Range("B10").Value() = $false

# Row 17 (4G+Cellular, Main Link section): no longer enabled
$ws.Range("B17").Value() = $false

# Update the active selection / view to G9 (and drop the A4 scroll position)
$null = $ws.Range("A1").Select()
$null = $ws.Range("G9").Select()
